# Applies the "Add files via upload" edit to the word list workbook:
#   1. Clears the existing AutoFilter criteria on column B (was filtering
#      for "Sb-105"), which un-hides every previously filtered-out row.
#   2. Inserts one new row (row 16) for the "exfoliqting" typo variant of
#      "exfoliating", mapped to the same "Stunin" key as the row above it.
#   3. Re-applies the AutoFilter over the now-larger A1:B50 range (with no
#      active filter criteria) and keeps the _FilterDatabase defined name
#      and the selected cell in sync with that.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the current filter criteria / unhide filtered rows ---------
$ws.ShowAllData()

# --- 2. Insert the new "exfoliqting" / "Stunin" row at row 16 ------------
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "exfoliqting"
$ws.Range("B16").Value = "Stunin"

# --- 3. Re-apply AutoFilter over the expanded range -----------------------
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:B50").AutoFilter()

# Keep the hidden _xlnm._FilterDatabase defined name in sync with the
# AutoFilter's new range (Sheet1!$A$1:$B$50).
$names = $wb.Names
$filterDbName = $names.Item("Sheet1!_FilterDatabase")
$filterDbName.RefersTo = "=Sheet1!`$A`$1:`$B`$50"

# Match the saved selection in the target workbook.
$ws.Range("A2").Select()
